# Atualização automática: 2025-08-24 09:00:44
#
# - Row 16 / Row 17: re-detection refined the bounding box / image file for
#   the same fly, updating D16/D17 (image filename) and I16 (coords), and
#   I17/J17 (coords + confidence).
# - A brand-new detection (row 21) is appended to the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 16 updates ----
$ws.Cells.Item(16, 4).Value = "image_20250807111728_ppp0.jpg"

$ws.Cells.Item(16, 9).Value = "'642,530,686,574"
$ws.Cells.Item(16, 9).Style = "Normal"

# ---- Row 17 updates ----
$ws.Cells.Item(17, 4).Value = "image_20250807111728_ppp0.jpg"

$ws.Cells.Item(17, 9).Value = "'794,481,830,525"
$ws.Cells.Item(17, 9).Style = "Normal"

$ws.Cells.Item(17, 10).Value = "'0.70"
$ws.Cells.Item(17, 10).Style = "Normal"

# ---- New row 21 ----
$ws.Cells.Item(21, 1).Value = "a2ea21b8-7dce-4e6a-be35-4edaddca5896"
$ws.Cells.Item(21, 2).Value = "mosca"

$ws.Cells.Item(21, 3).Value = 45893
$ws.Cells.Item(21, 3).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(21, 4).Value = "image_20250824092407_ppp0.jpg"
$ws.Cells.Item(21, 5).Value = "PLACA_20250717165933"
$ws.Cells.Item(21, 6).Value = "Beja"
$ws.Cells.Item(21, 7).Value = 38.02035
$ws.Cells.Item(21, 8).Value = -7.94715

$ws.Cells.Item(21, 9).Value = "'1002,789,1039,825"
$ws.Cells.Item(21, 9).Style = "Normal"

$ws.Cells.Item(21, 10).Value = "'0.64"
$ws.Cells.Item(21, 10).Style = "Normal"
